$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69, shifting existing rows 69-142 down to 70-143
$ws.Rows.Item(69).EntireRow.Insert()

# Populate the newly inserted row 69 with data
$ws.Cells.Item(69, 1).Value = 1
$ws.Cells.Item(69, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(69, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(69, 4).Value = 44966
$ws.Cells.Item(69, 5).Value = 15
$ws.Cells.Item(69, 6).Value = "Fruta"
$ws.Cells.Item(69, 7).Value = 100102
$ws.Cells.Item(69, 8).Value = "Cítricos"
$ws.Cells.Item(69, 9).Value = 100102004
$ws.Cells.Item(69, 10).Value = "Mandarina"
$ws.Cells.Item(69, 11).Value = "Murcott"
$ws.Cells.Item(69, 12).Value = "Tercera"
$ws.Cells.Item(69, 13).Value = 270
$ws.Cells.Item(69, 14).Value = 18000
$ws.Cells.Item(69, 15).Value = 20000
$ws.Cells.Item(69, 16).Value = 19000
$ws.Cells.Item(69, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(69, 18).Value = "Región Metropolitana"
$ws.Cells.Item(69, 19).Value = 950
$ws.Cells.Item(69, 20).Value = 20

Write-Output "done"
